# Update the workbook to add data for 2022-09-04:
# - rename sheet and update the "August 2022 (through August NN)" header
#   from "26" to "27"
# - update/insert the carjacking counts that changed in the B (August 2022
#   partial-month) column and scattered historical columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab name: "Through 2022-08-26" -> "Through 2022-08-27"
$ws.Name = "Through 2022-08-27"

# Column header text in row 1 (shared string), B1:
# "August 2022 (through August 26)" -> "August 2022 (through August 27)"
$ws.Range("B1").Value = "August 2022 (through August 27)"

# Updated / newly-populated data cells
$ws.Range("AP2").Value = 6
$ws.Range("R4").Value = 4
$ws.Range("AH4").Value = 5
$ws.Range("AX5").Value = 8
$ws.Range("AP6").Value = 2
$ws.Range("J7").Value = 6
$ws.Range("R7").Value = 7
$ws.Range("R9").Value = 6
$ws.Range("AX14").Value = 3
$ws.Range("B15").Value = 9
$ws.Range("R15").Value = 4
$ws.Range("J20").Value = 1
$ws.Range("AH24").Value = 1
$ws.Range("J27").Value = 2
$ws.Range("R29").Value = 3
$ws.Range("AH29").Value = 1
$ws.Range("R39").Value = 4
$ws.Range("B40").Value = 1
$ws.Range("R46").Value = 4
$ws.Range("AX50").Value = 3
$ws.Range("B55").Value = 2
$ws.Range("AH58").Value = 1
$ws.Range("J60").Value = 2
$ws.Range("B63").Value = 1
$ws.Range("B74").Value = 2
$ws.Range("BF92").Value = 2
$ws.Range("AH96").Value = 2
